$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad'
$ws.Range('G3').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range('G4').Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid'
$ws.Range('G5').Value = 'Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range('G6').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad'
$ws.Range('G7').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad'
$ws.Range('G8').Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator'
$ws.Range('G9').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Manar Montaser, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy'
$ws.Range('G10').Value = 'Dr. Gehan Adel, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range('G11').Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda'
$ws.Range('G13').Value = 'Dr. Omnia Mohammad, Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G14').Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G17').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa'
$ws.Range('G22').Value = 'Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy'
$ws.Range('G23').Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range('G24').Value = 'Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Salma Hassan, Dr. Aya Emad, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range('G25').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Aya Emad, Dr. Remon'
$ws.Range('G27').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Wafaa Ebida'
$ws.Range('G28').Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Remon, Dr. Aya Hanafy'
$ws.Range('G29').Value = 'Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon'
$ws.Range('G30').Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid'
$ws.Range('G31').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy'
$ws.Range('G32').Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid'
$ws.Range('G33').Value = 'Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Nesma, Dr. Mohammad El-Tanany, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Servinaz Sayed Mohammad'
$ws.Range('G34').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Servinaz Sayed Mohammad'
$ws.Range('G35').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Servinaz Sayed Mohammad'
$ws.Range('G36').Value = 'Dr. Asmaa Reda, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Manar Montaser, Dr. Majorelle Magdy, Administrator'
$ws.Range('G37').Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Manar Montaser, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Majorelle Magdy'
$ws.Range('G38').Value = 'Dr. Gehan Adel, Dr. Sara Wael, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Range('G39').Value = 'Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda'
$ws.Range('G41').Value = 'Dr. Omnia Mohammad, Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G42').Value = 'Dr. Safa Hany, Dr. Shimaa Ashraf'
$ws.Range('G45').Value = 'Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa'
$ws.Range('G50').Value = 'Dr. Amr Saeed, Dr. Nancy Abd Al-Shafy'
$ws.Range('G51').Value = 'Dr. Nourham Mostafa, Dr. Hana Amr'
$ws.Range('G52').Value = 'Dr. Monica, Dr. Remon, Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Neveen Nashaat, Dr. Maryam Ashraf, Dr. Marina Atef, Dr. Salma Hassan, Dr. Aya Emad, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range('G53').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Marina Atef, Dr. Aya Emad, Dr. Remon'
$ws.Range('G55').Value = 'Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Eman Mohammad Al, Dr. Neveen Nashaat, Dr. Yasmin, Dr. Salma Hassan, Dr. Remon, Dr. Wafaa Ebida'
$ws.Range('G56').Value = 'Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Nardine, Dr. Remon, Dr. Aya Hanafy'
$ws.Range('G57').Value = 'Dr. Monica, Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Remon'
